# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.873.03'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.636.56'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.73%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5023'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('E7').Value = '  -1.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2563'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06365'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.48'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07763'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.648.17'
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.256'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').Value = '1.862.01'
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5403'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').Value = '0.0₅7850'
$ws.Range('E16').Value = '  -0.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.54'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.69%  '
$ws.Range('D18').Value = '25.892.61'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '197.84'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.368'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.894'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.960'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.874'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '139.73'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1137'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.826'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.65'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.47%  '
$ws.Range('E30').Value = '  -0.55%  '
$ws.Range('E31').Value = '  -4.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.249'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.178'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.527'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.354'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.8857'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.600'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5526'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.92%  '
$ws.Range('D39').Value = '1.124.90'
$ws.Range('E39').Value = '  -0.87%  '
$ws.Range('E40').Value = '  -0.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9993'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.661'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8096'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.24'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.50%  '
$ws.Range('E45').Value = '  +9.07%  '
$ws.Range('D46').Value = '1.771.41'
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4518'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.28%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.006'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.60%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.16'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05070'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.90%  '
$ws.Range('E51').Value = '  -0.79%  '
